$wb = $excel.ActiveWorkbook

# Add the new "WMT_Extract_SA" worksheet after the last existing sheet (T2A)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "WMT_Extract_SA"

# Populate the header row
$newSheet.Range("A1").Value = "Case_Ref_No"
$newSheet.Range("B1").Value = "Tier_Code"
$newSheet.Range("C1").Value = "Team_Code"
$newSheet.Range("D1").Value = "OM_Grade_Code"
$newSheet.Range("E1").Value = "OM_Key"
$newSheet.Range("F1").Value = "Location"
$newSheet.Range("G1").Value = "Disposal_Type_Desc"
$newSheet.Range("H1").Value = "Disposal_Type_Code"
$newSheet.Range("I1").Value = "Standalone_Order"

# Copy the header formatting (bold white Arial on purple fill, text format)
# from an existing styled header cell, then recolour the border.
$srcStyle = $wb.Worksheets.Item("Court_Reports").Range("F1")
$srcStyle.Copy()
$newSheet.Range("A1:I1").PasteSpecial(-4122)
$newSheet.Range("A1:I1").Borders.Color = 16764108

# B1:I1 don't have a left border (it's supplied by the previous cell's right border)
foreach ($col in @("B", "C", "D", "E", "F", "G", "H", "I")) {
    $newSheet.Range($col + "1").Borders.Item(7).LineStyle = -4142
}

# Make the new sheet the active / selected tab, matching the source workbook
$newSheet.Activate()
$newSheet.Range("A1:I1").Select()
